$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New data for rows 2-9 (row 10 removed entirely)
# Columns: A (iteration idx), B (sample name string - stays same), C, D, E (numeric), L (S label string - shifts down by one)

$data = @(
    @{A=0; B="0m_1";  C=99.32416523349828; D=68.28859604845688; E=10.2277725699452;  L="S1"},
    @{A=1; B="8m_1";  C=99.37021207967101; D=68.51697066256951; E=10.2433500050944;  L="S2"},
    @{A=2; B="32m_1"; C=98.7402667566686;  D=66.60510628094494; E=13.57176929140085; L="S3"},
    @{A=3; B="1h_1";  C=99.31888229250808; D=68.24273052095465; E=9.346941192727559; L="S4"},
    @{A=4; B="4h_1";  C=98.83333010442296; D=68.06867810909776; E=11.17891278609345; L="S5"},
    @{A=5; B="8h_1";  C=98.84969293180292; D=68.19941794255345; E=9.39398540638576;  L="S6"},
    @{A=6; B="16h_1"; C=99.15328876853674; D=68.30072837970894; E=15.52677299785912; L="S7"},
    @{A=7; B="40h_1"; C=98.64836954674628; D=66.08044396034251; E=14.15203097545656; L="S8"}
)

$rowIdx = 2
foreach ($row in $data) {
    $ws.Cells.Item($rowIdx, 1).Value = $row.A
    $ws.Cells.Item($rowIdx, 2).Value = $row.B
    $ws.Cells.Item($rowIdx, 3).Value = $row.C
    $ws.Cells.Item($rowIdx, 4).Value = $row.D
    $ws.Cells.Item($rowIdx, 5).Value = $row.E
    $ws.Cells.Item($rowIdx, 12).Value = $row.L
    $rowIdx++
}

# Remove the old row 10 entirely (shift rows up / delete its content)
$ws.Rows.Item(10).Delete()
